# New penyata template format:
# Reset the Merit/Demerit figures for the homeroom account statement back to
# zero (blank template), and replace the four specific competition-name
# labels with the generic placeholder "Nama Pertandingan" repeated down the
# "Penyertaan Pertandingan" section.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reset Merit / Demerit amounts to 0 in the first few sections ---
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 0
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0

$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("D24").Value = 0
$ws.Range("E24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0

# --- "Penyertaan Pertandingan" section: rename each row label to the
#     generic "Nama Pertandingan" placeholder and clear the figures ---
$ws.Range("C34").Value = "Nama Pertandingan"
$ws.Range("D34").Value = 0

$ws.Range("C35").Value = "Nama Pertandingan"
$ws.Range("E35").Value = 0

$ws.Range("C36").Value = "Nama Pertandingan"
$ws.Range("D36").Value = 0

$ws.Range("C37").Value = "Nama Pertandingan"

$ws.Range("C38").Value = "Nama Pertandingan"

$ws.Range("C39").Value = "Nama Pertandingan"

$ws.Range("C40").Value = "Nama Pertandingan"
